$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the "RESPONSIBLE" column (G) entirely - shifts everything to its
# right one column to the left, matching the modelo de inventario rework.
$ws.Columns("G:G").Delete()

# Update the view: top-left visible cell, zoom and active selection.
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("G1").Select()
$excel.ActiveWindow.ScrollColumn = $ws.Range("G1").Column
$ws.Range("Q5").Select()
